$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells in columns D and E (Price / Volume) keep their original
# text representation instead of being auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.615.64"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.840.95"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.26"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4246"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3611"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07292"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8754"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.58"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.828.61"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.66%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.495"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06950"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.19"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008920"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.34"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.517.96"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.966"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.35"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.049.37"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.996"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.61"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "119.52"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.206"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.865"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08876"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7611"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.938"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.486"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.96%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05418"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.811"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1659"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5056"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.524"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.363"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06549"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "105.98"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.32"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4630"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.635"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.48%  "
